$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2061328790459966
$ws.Range("C2").Value = 0.5264054514480409
$ws.Range("J2").Value = 0.02214650766609881
$ws.Range("P2").Value = 0.1328790459965928
$ws.Range("S2").Value = 0.1124361158432709
$ws.Range("B3").Value = 0.006134969325153374
$ws.Range("C3").Value = 0.01840490797546012
$ws.Range("J3").Value = 0.02760736196319018
$ws.Range("P3").Value = 0.7361963190184049
$ws.Range("S3").Value = 0.2116564417177914
$ws.Range("J4").Value = 0.06578947368421052
$ws.Range("P4").Value = 0.6447368421052632
$ws.Range("S4").Value = 0.2894736842105263
$ws.Range("B6").Value = 0.05286343612334802
$ws.Range("D6").Value = 0.00881057268722467
$ws.Range("F6").Value = 0.07268722466960352
$ws.Range("J6").Value = 0.3325991189427313
$ws.Range("O6").Value = 0.02422907488986784
$ws.Range("Q6").Value = 0.1431718061674009
$ws.Range("R6").Value = 0.07488986784140969
$ws.Range("S6").Value = 0.2907488986784141
$ws.Range("B7").Value = 0.1113636363636364
$ws.Range("D7").Value = 0.01363636363636364
$ws.Range("E7").Value = 0.004545454545454545
$ws.Range("F7").Value = 0.06818181818181818
$ws.Range("J7").Value = 0.1318181818181818
$ws.Range("O7").Value = 0.01818181818181818
$ws.Range("Q7").Value = 0.1636363636363636
$ws.Range("R7").Value = 0.09090909090909091
$ws.Range("S7").Value = 0.3977272727272727
$ws.Range("B8").Value = 0.09834123222748815
$ws.Range("D8").Value = 0.01658767772511848
$ws.Range("F8").Value = 0.06990521327014218
$ws.Range("J8").Value = 0.1054502369668247
$ws.Range("O8").Value = 0.02488151658767773
$ws.Range("Q8").Value = 0.1670616113744076
$ws.Range("R8").Value = 0.1007109004739336
$ws.Range("S8").Value = 0.4170616113744076
$ws.Range("B9").Value = 0.1253481894150418
$ws.Range("D9").Value = 0.01949860724233983
$ws.Range("F9").Value = 0.08635097493036212
$ws.Range("J9").Value = 0.1420612813370473
$ws.Range("O9").Value = 0.01114206128133705
$ws.Range("Q9").Value = 0.1727019498607242
$ws.Range("R9").Value = 0.08356545961002786
$ws.Range("S9").Value = 0.3593314763231198
$ws.Range("B10").Value = 0.1030502885408079
$ws.Range("D10").Value = 0.02019785655399835
$ws.Range("E10").Value = 0.0008244023083264633
$ws.Range("F10").Value = 0.05935696619950536
$ws.Range("J10").Value = 0.1084089035449299
$ws.Range("O10").Value = 0.01483924154987634
$ws.Range("Q10").Value = 0.2328936521022259
$ws.Range("R10").Value = 0.08615004122011542
$ws.Range("S10").Value = 0.3742786479802143
$ws.Range("G11").Value = 0.1597717546362339
$ws.Range("J11").Value = 0.08701854493580599
$ws.Range("K11").Value = 0.2082738944365193
$ws.Range("L11").Value = 0.5306704707560628
$ws.Range("S11").Value = 0.01426533523537803
$ws.Range("G12").Value = 0.741514360313316
$ws.Range("J12").Value = 0.2010443864229765
$ws.Range("K12").Value = 0.01044386422976501
$ws.Range("L12").Value = 0.02088772845953003
$ws.Range("S12").Value = 0.02610966057441253
$ws.Range("G13").Value = 0.6236559139784946
$ws.Range("J13").Value = 0.3763440860215054
$ws.Range("F15").Value = 0.03178484107579462
$ws.Range("H15").Value = 0.136919315403423
$ws.Range("I15").Value = 0.06112469437652811
$ws.Range("J15").Value = 0.352078239608802
$ws.Range("K15").Value = 0.07090464547677261
$ws.Range("M15").Value = 0.007334963325183374
$ws.Range("O15").Value = 0.04645476772616137
$ws.Range("S15").Value = 0.293398533007335
$ws.Range("F16").Value = 0.01685393258426966
$ws.Range("H16").Value = 0.1629213483146068
$ws.Range("I16").Value = 0.08426966292134831
$ws.Range("J16").Value = 0.4157303370786517
$ws.Range("K16").Value = 0.1432584269662921
$ws.Range("M16").Value = 0.01685393258426966
$ws.Range("N16").Value = 0.002808988764044944
$ws.Range("O16").Value = 0.0449438202247191
$ws.Range("S16").Value = 0.1123595505617977
$ws.Range("F17").Value = 0.02669632925472748
$ws.Range("H17").Value = 0.167964404894327
$ws.Range("I17").Value = 0.1012235817575083
$ws.Range("J17").Value = 0.3904338153503893
$ws.Range("K17").Value = 0.1112347052280311
$ws.Range("M17").Value = 0.02224694104560623
$ws.Range("N17").Value = 0.001112347052280311
$ws.Range("O17").Value = 0.04338153503893215
$ws.Range("S17").Value = 0.135706340378198
$ws.Range("F18").Value = 0.03535353535353535
$ws.Range("H18").Value = 0.1464646464646465
$ws.Range("I18").Value = 0.07323232323232323
$ws.Range("J18").Value = 0.3838383838383838
$ws.Range("K18").Value = 0.143939393939394
$ws.Range("M18").Value = 0.01262626262626263
$ws.Range("O18").Value = 0.08080808080808081
$ws.Range("S18").Value = 0.1237373737373737
$ws.Range("F19").Value = 0.02279202279202279
$ws.Range("H19").Value = 0.2157102157102157
$ws.Range("I19").Value = 0.07651607651607652
$ws.Range("J19").Value = 0.3508343508343508
$ws.Range("K19").Value = 0.1245421245421245
$ws.Range("M19").Value = 0.02604802604802605
$ws.Range("N19").Value = 0.000814000814000814
$ws.Range("O19").Value = 0.06552706552706553
$ws.Range("S19").Value = 0.1172161172161172
